$wb = $excel.ActiveWorkbook

# OFF sheet - Week 16 logged values (row 3)
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 456
$wsOff.Range("C3").Value = 311
$wsOff.Range("D3").Value = 93
$wsOff.Range("E3").Value = 39
$wsOff.Range("G3").Value = 7

# DEF sheet - Week 16 logged values (row 3)
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 470
$wsDef.Range("C3").Value = 324
$wsDef.Range("D3").Value = 114
$wsDef.Range("E3").Value = 54
